$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "289.39"
Set-TextValue 2 5 "-4.15%"
Set-TextValue 3 4 "30.78"
Set-TextValue 3 5 "-4.19%"
Set-TextValue 4 4 "4.878"
Set-TextValue 4 5 "-2.13%"
Set-TextValue 5 4 "0.07153"
Set-TextValue 5 5 "-9.49%"
Set-TextValue 6 4 "1.865"
Set-TextValue 6 5 "-11.96%"
Set-TextValue 7 4 "7.639"
Set-TextValue 8 4 "3.724"
Set-TextValue 8 5 "-1.94%"
Set-TextValue 9 4 "0.8976"
Set-TextValue 9 5 "-3.10%"
Set-TextValue 10 4 "0.1649"
Set-TextValue 10 5 "-5.92%"
Set-TextValue 11 4 "0.07499"
Set-TextValue 11 5 "-5.48%"
Set-TextValue 12 4 "0.08106"
Set-TextValue 12 5 "-6.27%"
Set-TextValue 13 4 "0.02991"
Set-TextValue 13 5 "-4.62%"
Set-TextValue 14 4 "0.09983"
Set-TextValue 14 5 "-0.27%"
Set-TextValue 15 4 "0.001494"
Set-TextValue 15 5 "-1.95%"
Set-TextValue 16 4 "0.005816"
Set-TextValue 16 5 "0.88%"
Set-TextValue 18 4 "3.461"
Set-TextValue 19 4 "2.108"
Set-TextValue 19 5 "-7.39%"
Set-TextValue 20 4 "0.3277"
Set-TextValue 20 5 "-0.32%"
Set-TextValue 21 5 "-0.74%"
Set-TextValue 22 4 "4.305"
Set-TextValue 22 5 "0.85%"
Set-TextValue 23 4 "0.2002"
Set-TextValue 23 5 "11.77%"
Set-TextValue 24 4 "0.04479"
Set-TextValue 24 5 "-2.81%"
Set-TextValue 25 4 "0.001214"
Set-TextValue 25 5 "-1.90%"
Set-TextValue 26 4 "0.004657"
Set-TextValue 26 5 "4.33%"
Set-TextValue 27 4 "0.0001252"
Set-TextValue 27 5 "0.12%"
Set-TextValue 39 4 "0.01644"
Set-TextValue 39 5 "-4.51%"
Set-TextValue 40 4 "0.04343"
Set-TextValue 40 5 "-9.23%"
Set-TextValue 41 4 "0.007311"
Set-TextValue 41 5 "-1.63%"
Set-TextValue 42 4 "0.1303"
Set-TextValue 42 5 "-4.17%"
Set-TextValue 43 4 "0.002007"
Set-TextValue 43 5 "-16.03%"
Set-TextValue 44 4 "0.01023"
Set-TextValue 44 5 "-0.02%"
Set-TextValue 45 4 "0.00005814"
Set-TextValue 45 5 "-3.02%"
Set-TextValue 46 4 "0.00000000751"
Set-TextValue 46 5 "0.13%"
Set-TextValue 47 4 "2.202"
Set-TextValue 47 5 "167.40%"
Set-TextValue 48 5 "-11.45%"
Set-TextValue 49 4 "0.00002103"
Set-TextValue 49 5 "0.13%"
Set-TextValue 50 4 "0.0002003"
Set-TextValue 50 5 "0.13%"
